$d = $word.ActiveDocument

# The paragraph currently contains three runs:
#   <id>      (Courier New, color 7f6000, sz 18)
#   p152r_1   (color 000000)
#   </id>     (Courier New, color 7f6000, sz 18)
# They must become a single run "<id>p152r_1</id>" using the formatting
# of the first ("<id>") run.

$rngFull = $d.Content
$found = $rngFull.Find.Execute("<id>p152r_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $idStart = $rngFull.Start
    $idEnd = $rngFull.End

    # Keep the first run ("<id>") untouched; delete everything after it
    # within the match (i.e. "p152r_1</id>").
    $restRange = $d.Range($idStart + 4, $idEnd)
    $restRange.Delete()

    # Re-insert the remaining text right after "<id>". Because the
    # collapsed range sits at the end of the "<id>" run, it inherits that
    # run's formatting (Courier New, color 7f6000, size 18) automatically,
    # so the inserted text merges into the same run instead of creating a
    # new one.
    $restRange.InsertAfter("p152r_1</id>")
}
